$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 69333.336
$ws.Range("J52").Value = 69333.336
$ws.Range("L52").Value = 208000.008
$ws.Range("N52").Value = -208320.008
$ws.Range("H100").Value = 20001518
$ws.Range("I100").Value = 20001518
$ws.Range("K100").Value = 20001518
$ws.Range("M100").Value = -20000977
$ws.Range("H129").Value = 801.52527
$ws.Range("J129").Value = 841.4176
$ws.Range("L129").Value = 2524.2528
$ws.Range("N129").Value = -12524.2528
$ws.Range("H137").Value = 1324803.4
$ws.Range("I137").Value = 2269047.8
$ws.Range("J137").Value = 2861.2666
$ws.Range("K137").Value = 6807143.399999999
$ws.Range("L137").Value = 8583.799800000001
$ws.Range("M137").Value = -6804593.399999999
$ws.Range("N137").Value = -13683.7998
$ws.Range("H138").Value = 4930.31
$ws.Range("I138").Value = 807.8095
$ws.Range("J138").Value = 6026.1646
$ws.Range("K138").Value = 2423.4285
$ws.Range("L138").Value = 18078.4938
$ws.Range("M138").Value = 2716.5715
$ws.Range("N138").Value = -28358.4938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4826.769
$ws.Range("I32").Value = 4212.0215
$ws.Range("K32").Value = 4212.0215
$ws.Range("M32").Value = -3925.0215
$ws.Range("H45").Value = 2424.625
$ws.Range("I45").Value = 2720.6
$ws.Range("J45").Value = 1931.3334
$ws.Range("K45").Value = 2720.6
$ws.Range("L45").Value = 1931.3334
$ws.Range("M45").Value = -2343.6
$ws.Range("N45").Value = -2685.3334
$ws.Range("H61").Value = 1525.2858
$ws.Range("I61").Value = 1454.5834
$ws.Range("J61").Value = 1949.5
$ws.Range("K61").Value = 1454.5834
$ws.Range("L61").Value = 1949.5
$ws.Range("M61").Value = -1242.5834
$ws.Range("N61").Value = -2373.5
$ws.Range("H74").Value = 4168.3438
$ws.Range("I74").Value = 4670.905
$ws.Range("J74").Value = 3208.9092
$ws.Range("K74").Value = 4670.905
$ws.Range("L74").Value = 3208.9092
$ws.Range("M74").Value = -3796.905
$ws.Range("N74").Value = -4956.9092
$ws.Range("H77").Value = 4168.3438
$ws.Range("I77").Value = 4670.905
$ws.Range("J77").Value = 3208.9092
$ws.Range("K77").Value = 23354.525
$ws.Range("L77").Value = 16044.546
$ws.Range("M77").Value = -18986.525
$ws.Range("N77").Value = -24780.546
$ws.Range("H102").Value = 990
$ws.Range("I102").Value = 990
$ws.Range("K102").Value = 990
$ws.Range("M102").Value = 632
$ws.Range("H136").Value = 1525.2858
$ws.Range("I136").Value = 1454.5834
$ws.Range("J136").Value = 1949.5
$ws.Range("K136").Value = 4363.7502
$ws.Range("L136").Value = 5848.5
$ws.Range("M136").Value = -1813.7502
$ws.Range("N136").Value = -10948.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 718.52
$ws.Range("I94").Value = 607.087
$ws.Range("K94").Value = 607.087
$ws.Range("M94").Value = -156.087
$ws.Range("H103").Value = 30442.666
$ws.Range("J103").Value = 30442.666
$ws.Range("L103").Value = 30442.666
$ws.Range("N103").Value = -32786.666
$ws.Range("H105").Value = 4976734.5
$ws.Range("I105").Value = 5209943
$ws.Range("K105").Value = 5209943
$ws.Range("M105").Value = -5208196

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H31").Value = 2150.3914
$ws.Range("I31").Value = 1234.6842
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 1234.6842
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -939.6841999999999
$ws.Range("N31").Value = -7090
$ws.Range("H34").Value = 2150.3914
$ws.Range("I34").Value = 1234.6842
$ws.Range("J34").Value = 6500
$ws.Range("K34").Value = 1234.6842
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = -1032.6842
$ws.Range("N34").Value = -6904
$ws.Range("H122").Value = 2074.647
$ws.Range("I122").Value = 1077.1
$ws.Range("J122").Value = 3499.7144
$ws.Range("K122").Value = 3231.3
$ws.Range("L122").Value = 10499.1432
$ws.Range("M122").Value = -781.2999999999997
$ws.Range("N122").Value = -15399.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2277.8616
$ws.Range("J121").Value = 2344.4602
$ws.Range("L121").Value = 7033.3806
$ws.Range("N121").Value = -9653.3806
$ws.Range("H129").Value = 2373.6
$ws.Range("I129").Value = 2334.8333
$ws.Range("J129").Value = 2431.75
$ws.Range("K129").Value = 7004.499899999999
$ws.Range("L129").Value = 7295.25
$ws.Range("M129").Value = -2004.499899999999
$ws.Range("N129").Value = -17295.25
$ws.Range("H140").Value = 2416.8096
$ws.Range("I140").Value = 2416.8096
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 7250.4288
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2070.4288
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 29998
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 29998
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 29998
$ws.Range("N4").Value = -30222
$ws.Range("M4").ClearContents()
$ws.Range("H80").Value = 50002420
$ws.Range("I80").Value = 62502276
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 62502276
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -62501278
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 50002420
$ws.Range("I83").Value = 62502276
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 312511380
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -312506388
$ws.Range("N83").Value = -24984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7415
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 8581.25
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 25743.75
$ws.Range("M122").Value = -5800
$ws.Range("N122").Value = -30643.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 38333.332
$ws.Range("I62").Value = 8000
$ws.Range("J62").Value = 53500
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 53500
$ws.Range("M62").Value = -7376
$ws.Range("N62").Value = -54748
$ws.Range("H65").Value = 38333.332
$ws.Range("I65").Value = 8000
$ws.Range("J65").Value = 53500
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 267500
$ws.Range("M65").Value = -36880
$ws.Range("N65").Value = -273740
$ws.Range("H96").Value = 168417200
$ws.Range("I96").Value = 252625250
$ws.Range("J96").Value = 1099.5
$ws.Range("K96").Value = 252625250
$ws.Range("L96").Value = 1099.5
$ws.Range("M96").Value = -252623877
$ws.Range("N96").Value = -3845.5
$ws.Range("H114").Value = 27296
$ws.Range("J114").Value = 27296
$ws.Range("L114").Value = 27296
$ws.Range("N114").Value = -35974
$ws.Range("H122").Value = 6144.8887
$ws.Range("I122").Value = 3384
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 10152
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -7702
$ws.Range("N122").Value = -39900.001
